$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-11: C = 11, D = row-1 (1..10), remove F/G values
# Then add rows 12-29 with B=1, C=11, D=11..28, E=TRUE

for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
    $ws.Cells.Item($r, 3).Value = 11
    $ws.Cells.Item($r, 4).Value = ($r - 1)
    $ws.Cells.Item($r, 5).Value = $true
}

# Clear F and G columns for rows 2-11 (previously had 0 values, now removed)
$ws.Range("F2:G11").ClearContents()

# Set column widths to match bestFit columns from the diff
$ws.Columns.Item(1).ColumnWidth = 2.7109375
$ws.Columns.Item(2).ColumnWidth = 8.7109375
$ws.Columns.Item(3).ColumnWidth = 15.28515625
$ws.Columns.Item(4).ColumnWidth = 7.42578125
$ws.Columns.Item(5).ColumnWidth = 6.42578125
$ws.Columns.Item(6).ColumnWidth = 11.42578125
$ws.Columns.Item(7).ColumnWidth = 12.140625

# Set the active selection cell to I23 to match the diff
$ws.Range("I23").Select()
